$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Invoice Date / Due Date block (rows 12-13) ---
# Bump the invoice date and due date forward by one day, and give the
# label cells (A:D) the same date number-format already used by the
# date values in column E.
$ws.Range("A12:D13").NumberFormat = $ws.Range("E12").NumberFormat
$ws.Range("E12").Value = 44367
$ws.Range("E13").Value = 44367

# --- Footer ---
# The "Thank You" / "Terms & Instructions" notes used to live in their
# own standalone rows below the totals block; move them up into the
# (until now blank) column A of the Subtotal/Discounts rows instead.
$ws.Range("A28").Value = "Thank You for your business!"
$ws.Range("A29").Value = "Terms & Instructions"

# The old placeholder B/C cells for the totals rows are no longer needed.
$ws.Range("B28:C31").ClearContents()

# These rows were merged across A:C purely to hold the (now-removed)
# blank placeholder text, so undo the merges.
$ws.Range("A28:C28").UnMerge()
$ws.Range("A29:C29").UnMerge()
$ws.Range("A30:C30").UnMerge()
$ws.Range("A31:C31").UnMerge()

# The standalone footer rows are now redundant since their text moved
# into rows 28-29 above, so drop them.
$ws.Rows.Item(32).EntireRow.Delete()
$ws.Rows.Item(32).EntireRow.Delete()
